# Historical year electricity calibration
# Set the Boolean "new plant allowed" flags for years 2021-2023 (columns B, C, D)
# to 0 for all generator technologies except the CCS variants (rows 19-22),
# which were already 0 and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BBNPPTY")

$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,23,24,25)
foreach ($r in $rows) {
    $ws.Range("B${r}:D${r}").Value = 0
}

# Make BBNPPTY the active sheet / tab, matching the author's last saved view,
# with the selection left on G27.
$ws.Select()
$ws.Range("G27").Select()
